# Add two new columns, I ("I0") and J ("IF"), to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells (row 1) -------------------------------------------------
# Copy the formatting of the existing header cell H1 (bold, centered,
# bordered) onto the new header cells so the new columns match the look of
# the existing ones, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data cells (rows 2-58) ----------------------------------------------
$iValues = @(6,4,8,7,7,8,8,8,8,8,8,9,7,8,9,9,9,10,8,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,7,10,9,9,9,9,9,9,9,9,9,9,8,8,7,7,5,7,6,5,4,8,2)
$jValues = @(6,5,9,7,7,8,8,8,8,8,8,9,7,8,10,9,9,10,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,7,10,9,9,9,9,9,9,9,9,9,9,8,9,7,7,5,7,7,5,4,8,2)

$firstRow = 2
$lastRow = 58
$rowCount = $lastRow - $firstRow + 1

$iData = New-Object 'object[,]' $rowCount,1
$jData = New-Object 'object[,]' $rowCount,1
for ($idx = 0; $idx -lt $rowCount; $idx++) {
    $iData[$idx,0] = $iValues[$idx]
    $jData[$idx,0] = $jValues[$idx]
}

$ws.Range("I$firstRow`:I$lastRow").Value = $iData
$ws.Range("J$firstRow`:J$lastRow").Value = $jData
